$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text corrections
$ws.Range("B2").Value = "<them>"
$ws.Range("B17").Value = "<like>"

# Column C numeric corrections
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 7
$ws.Range("C5").Value = 9
$ws.Range("C6").Value = 8
$ws.Range("C7").Value = 4
$ws.Range("C9").Value = 5
$ws.Range("C12").Value = 7
$ws.Range("C13").Value = 6
$ws.Range("C14").Value = 7
$ws.Range("C15").Value = 8
$ws.Range("C16").Value = 8
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 5
